$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 869
$ws.Cells.Item(869, 1).Value = "SCLD Student Event"
$ws.Cells.Item(869, 2).Value = 42676
$ws.Cells.Item(869, 3).Value = "1800"
$ws.Cells.Item(869, 4).Value = "WC"
$ws.Cells.Item(869, 5).Value = "118"
$ws.Cells.Item(869, 6).Value = "INC000000739329"

# Row 870
$ws.Cells.Item(870, 1).Value = "SCLD Student Logout"
$ws.Cells.Item(870, 2).Value = 42676
$ws.Cells.Item(870, 3).Value = "2100"
$ws.Cells.Item(870, 4).Value = "WC"
$ws.Cells.Item(870, 5).Value = "118"
$ws.Cells.Item(870, 6).Value = "INC000000739329"

# Row 871
$ws.Cells.Item(871, 1).Value = "SCLD Student Event"
$ws.Cells.Item(871, 2).Value = 42676
$ws.Cells.Item(871, 3).Value = "1830"
$ws.Cells.Item(871, 4).Value = "FC"
$ws.Cells.Item(871, 5).Value = "104"
$ws.Cells.Item(871, 6).Value = "INC000000733582"

# Row 872
$ws.Cells.Item(872, 1).Value = "SCLD Student Logout"
$ws.Cells.Item(872, 2).Value = 42676
$ws.Cells.Item(872, 3).Value = "2030"
$ws.Cells.Item(872, 4).Value = "FC"
$ws.Cells.Item(872, 5).Value = "104"
$ws.Cells.Item(872, 6).Value = "INC000000733582"

# Row 877
$ws.Cells.Item(877, 1).Value = "Demo"
$ws.Cells.Item(877, 2).Value = 42677
$ws.Cells.Item(877, 3).Value = "1700"
$ws.Cells.Item(877, 4).Value = "SSB"
$ws.Cells.Item(877, 5).Value = "W141"
$ws.Cells.Item(877, 6).Value = "Client usin PC, podium mic and 2 neck mics (one built in, second in the back booth, plugged into mixing board)"

# Row 878
$ws.Cells.Item(878, 1).Value = "Operator"
$ws.Cells.Item(878, 2).Value = 42677
$ws.Cells.Item(878, 3).Value = "1715"
$ws.Cells.Item(878, 4).Value = "SSB"
$ws.Cells.Item(878, 5).Value = "W141"
$ws.Cells.Item(878, 6).Value = "Operate event / assist client between 5:15-5:45"

# Row 879
$ws.Cells.Item(879, 1).Value = "Setup Mic"
$ws.Cells.Item(879, 2).Value = 42677
$ws.Cells.Item(879, 3).Value = "1800"
$ws.Cells.Item(879, 4).Value = "DB"
$ws.Cells.Item(879, 5).Value = "2027"
$ws.Cells.Item(879, 6).Value = "Neck mic and small PA from DB 0003"

# Row 880
$ws.Cells.Item(880, 1).Value = "Pickup Mic"
$ws.Cells.Item(880, 2).Value = 42677
$ws.Cells.Item(880, 3).Value = "2100"
$ws.Cells.Item(880, 4).Value = "DB"
$ws.Cells.Item(880, 5).Value = "2027"
$ws.Cells.Item(880, 6).Value = "Return neck mic and small PA to DB 0003"

# Row 881
$ws.Cells.Item(881, 1).Value = "Demo"
$ws.Cells.Item(881, 2).Value = 42677
$ws.Cells.Item(881, 3).Value = "1830"
$ws.Cells.Item(881, 4).Value = "OSG"
$ws.Cells.Item(881, 5).Value = "1005"

# Row 882
$ws.Cells.Item(882, 1).Value = "Demo"
$ws.Cells.Item(882, 2).Value = 42677
$ws.Cells.Item(882, 3).Value = "1900"
$ws.Cells.Item(882, 4).Value = "SSB"
$ws.Cells.Item(882, 5).Value = "N108"
$ws.Cells.Item(882, 6).Value = "Using neck mic"

# Row 883
$ws.Cells.Item(883, 1).Value = "SCLD Student Event"
$ws.Cells.Item(883, 2).Value = 42676
$ws.Cells.Item(883, 3).Value = "1800"
$ws.Cells.Item(883, 4).Value = "CLH"
$ws.Cells.Item(883, 5).Value = "K"
$ws.Cells.Item(883, 6).Value = "INC000000740762"

# Row 884
$ws.Cells.Item(884, 1).Value = "SCLD Student Logout"
$ws.Cells.Item(884, 2).Value = 42676
$ws.Cells.Item(884, 3).Value = "2030"
$ws.Cells.Item(884, 4).Value = "CLH"
$ws.Cells.Item(884, 5).Value = "K"
$ws.Cells.Item(884, 6).Value = "INC000000740762"

# Row 888
$ws.Cells.Item(888, 1).Value = "Other"
$ws.Cells.Item(888, 2).Value = 42678
$ws.Cells.Item(888, 3).Value = "1700"
$ws.Cells.Item(888, 4).Value = "OSG"
$ws.Cells.Item(888, 5).Value = "1014"
$ws.Cells.Item(888, 6).Value = "Turn off mixer / leave in place everything for next day"

# Row 889
$ws.Cells.Item(889, 1).Value = "Pickup Mic"
$ws.Cells.Item(889, 2).Value = 42678
$ws.Cells.Item(889, 3).Value = "1600"
$ws.Cells.Item(889, 4).Value = "KT"
$ws.Cells.Item(889, 5).Value = "280 N "
$ws.Cells.Item(889, 6).Value = "Room in York Lanes - return 2 IR mics to KT 516 and place batteries in charger"

# Row 890
$ws.Cells.Item(890, 1).Value = "Pickup Mic"
$ws.Cells.Item(890, 2).Value = 42678
$ws.Cells.Item(890, 3).Value = "1730"
$ws.Cells.Item(890, 4).Value = "SSB"
$ws.Cells.Item(890, 5).Value = "W141"
$ws.Cells.Item(890, 6).Value = "Pick up one audience handheld mic  and stands and return to rear booth. Leave podium mic  and desk mics in place"

# Row 895
$ws.Cells.Item(895, 1).Value = "Demo"
$ws.Cells.Item(895, 2).Value = 42681
$ws.Cells.Item(895, 3).Value = "1900"
$ws.Cells.Item(895, 4).Value = "SSB"
$ws.Cells.Item(895, 5).Value = "S124"

# Rows whose long wrapped comment text forces a taller (30pt) row in Excel
$ws.Rows.Item(877).RowHeight = 30
$ws.Rows.Item(889).RowHeight = 30
$ws.Rows.Item(890).RowHeight = 30

# Update the view to match where the user ended up after typing the new entries
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 881
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E900").Select() | Out-Null

